# añadir autoevaluacion de Megane
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in self-evaluation values (Yo/D) that were previously missing
$ws.Range("D24").Value = 2

# Fill in values for participant E (E23, E24, E26:E29)
$ws.Range("E23").Value = 1
$ws.Range("E24").Value = 3

# Fill in values for participant F (F23, F24, F26:F29)
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 3

# Row 26 (indicator 16) values across D, E, F
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1

# Row 27 (indicator 17) values across D, E, F
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 3

# Row 28 (indicator 18) values across D, E, F
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = 3
$ws.Range("F28").Value = 3

# Row 29 (indicator 19) values across D, E, F
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 1

# Update the view to reflect the scrolled position and new selection
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws.Range("E24").Select()

$wb.Save()
